# The document contains four occurrences of a pattern where the literal
# text "<id>p038r_N</id>" is split across three separate runs:
#   run 1: "<id>"        (Courier New, color 7f6000, sz 18)
#   run 2: "p038r_N"     (color 000000)
#   run 3: "</id>"       (Courier New, color 7f6000, sz 18)
#
# The edit merges each trio of runs into a single run containing the
# full text "<id>p038r_N</id>", keeping the formatting of the first run
# (Courier New / color 7f6000 / sz 18). Using Find/Replace across the
# exact text span causes Word to replace the whole matched range with a
# single new run, which is exactly the collapse we need.

$d = $word.ActiveDocument

for ($i = 1; $i -le 4; $i++) {
    $needle = "<id>p038r_$i</id>"
    $found = $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
    if (-not $found) {
        Write-Host "WARNING: pattern not found for p038r_$i"
    }
}
